# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
#
# Both sheets list the same set of events (in the same order), but
# "全部类型" has one extra row inserted part-way down, so absolute row
# numbers differ between the two sheets. To make the update robust we
# walk column F top-to-bottom on each sheet and, for each still-unmatched
# (old -> new) pair, apply it to the first row whose current value equals
# the expected "old" value. Processing the pairs in the same order they
# appear in the source data (top to bottom) guarantees the correct row is
# updated even when two rows momentarily share the same old value.

$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Old = 558;   New = 561 },
    @{ Old = 1123;  New = 1124 },
    @{ Old = 64;    New = 65 },
    @{ Old = 51;    New = 55 },
    @{ Old = 1155;  New = 1156 },
    @{ Old = 16185; New = 16220 },
    @{ Old = 269;   New = 271 },
    @{ Old = 198;   New = 199 },
    @{ Old = 6324;  New = 6336 },
    @{ Old = 18;    New = 19 },
    @{ Old = 21;    New = 23 },
    @{ Old = 34;    New = 36 },
    @{ Old = 20;    New = 22 },
    @{ Old = 12;    New = 14 },
    @{ Old = 46;    New = 47 },
    @{ Old = 5038;  New = 5041 },
    @{ Old = 494;   New = 495 },
    @{ Old = 11275; New = 11288 },
    @{ Old = 143;   New = 145 },
    @{ Old = 196;   New = 197 },
    @{ Old = 3831;  New = 3832 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.UsedRange.Rows.Count

    foreach ($change in $changes) {
        for ($r = 2; $r -le $lastRow; $r++) {
            $cell = $ws.Cells.Item($r, 6)  # column F = "想去人数"
            if ($cell.Value2 -eq $change.Old) {
                $cell.Value2 = $change.New
                break
            }
        }
    }
}
